$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp title
$ws.Cells.Item(1,1).Value = "Datos actualizados a 28 de Marzo de 2020 a las 15:59"

# Row 4: Estados Unidos
$ws.Cells.Item(4,1).Value = "Estados Unidos"
$ws.Cells.Item(4,2).Value = 105019
$ws.Cells.Item(4,3).Value = 893
$ws.Cells.Item(4,4).Value = 2537
$ws.Cells.Item(4,5).Value = 100765
$ws.Cells.Item(4,6).Value = 2494
$ws.Cells.Item(4,7).Value = 21
$ws.Cells.Item(4,8).Value = 1717

# Row 20: Noruega
$ws.Cells.Item(20,1).Value = "Noruega"
$ws.Cells.Item(20,2).Value = 3972
$ws.Cells.Item(20,3).Value = 201
$ws.Cells.Item(20,4).Value = 7
$ws.Cells.Item(20,5).Value = 3945
$ws.Cells.Item(20,6).Value = 76
$ws.Cells.Item(20,7).Value = 1
$ws.Cells.Item(20,8).Value = 20

# Row 46: Singapur
$ws.Cells.Item(46,1).Value = "Singapur"
$ws.Cells.Item(46,2).Value = 802
$ws.Cells.Item(46,3).Value = 70
$ws.Cells.Item(46,4).Value = 198
$ws.Cells.Item(46,5).Value = 602
$ws.Cells.Item(46,6).Value = 19
$ws.Cells.Item(46,7).Value = 0
$ws.Cells.Item(46,8).Value = 2

# Row 47: Panama
$ws.Cells.Item(47,1).Value = "Panama"
$ws.Cells.Item(47,2).Value = 786
$ws.Cells.Item(47,3).Value = 0
$ws.Cells.Item(47,4).Value = 2
$ws.Cells.Item(47,5).Value = 770
$ws.Cells.Item(47,6).Value = 20
$ws.Cells.Item(47,7).Value = 0
$ws.Cells.Item(47,8).Value = 14

# Row 48: Republica Dominicana
$ws.Cells.Item(48,1).Value = "Republica Dominicana"
$ws.Cells.Item(48,2).Value = 719
$ws.Cells.Item(48,3).Value = 138
$ws.Cells.Item(48,4).Value = 3
$ws.Cells.Item(48,5).Value = 688
$ws.Cells.Item(48,6).Value = 0
$ws.Cells.Item(48,7).Value = 8
$ws.Cells.Item(48,8).Value = 28

# Row 49: Mexico
$ws.Cells.Item(49,1).Value = "Mexico"
$ws.Cells.Item(49,2).Value = 717
$ws.Cells.Item(49,3).Value = 132
$ws.Cells.Item(49,4).Value = 4
$ws.Cells.Item(49,5).Value = 701
$ws.Cells.Item(49,6).Value = 1
$ws.Cells.Item(49,7).Value = 4
$ws.Cells.Item(49,8).Value = 12

# Row 50: Crucero
$ws.Cells.Item(50,1).Value = "Crucero"
$ws.Cells.Item(50,2).Value = 712
$ws.Cells.Item(50,3).Value = 0
$ws.Cells.Item(50,4).Value = 597
$ws.Cells.Item(50,5).Value = 105
$ws.Cells.Item(50,6).Value = 15
$ws.Cells.Item(50,7).Value = 0
$ws.Cells.Item(50,8).Value = 10

# Row 51: Argentina
$ws.Cells.Item(51,1).Value = "Argentina"
$ws.Cells.Item(51,2).Value = 690
$ws.Cells.Item(51,3).Value = 101
$ws.Cells.Item(51,4).Value = 72
$ws.Cells.Item(51,5).Value = 601
$ws.Cells.Item(51,6).Value = 0
$ws.Cells.Item(51,7).Value = 4
$ws.Cells.Item(51,8).Value = 17

# Row 52: Eslovenia
$ws.Cells.Item(52,1).Value = "Eslovenia"
$ws.Cells.Item(52,2).Value = 684
$ws.Cells.Item(52,3).Value = 52
$ws.Cells.Item(52,4).Value = 10
$ws.Cells.Item(52,5).Value = 665
$ws.Cells.Item(52,6).Value = 25
$ws.Cells.Item(52,7).Value = 0
$ws.Cells.Item(52,8).Value = 9

# Row 53: Serbia
$ws.Cells.Item(53,1).Value = "Serbia"
$ws.Cells.Item(53,2).Value = 659
$ws.Cells.Item(53,3).Value = 131
$ws.Cells.Item(53,4).Value = 42
$ws.Cells.Item(53,5).Value = 607
$ws.Cells.Item(53,6).Value = 25
$ws.Cells.Item(53,7).Value = 2
$ws.Cells.Item(53,8).Value = 10

# Row 54: Estonia
$ws.Cells.Item(54,1).Value = "Estonia"
$ws.Cells.Item(54,2).Value = 645
$ws.Cells.Item(54,3).Value = 70
$ws.Cells.Item(54,4).Value = 20
$ws.Cells.Item(54,5).Value = 624
$ws.Cells.Item(54,6).Value = 10
$ws.Cells.Item(54,7).Value = 0
$ws.Cells.Item(54,8).Value = 1

# Row 55: Peru
$ws.Cells.Item(55,1).Value = "Peru"
$ws.Cells.Item(55,2).Value = 635
$ws.Cells.Item(55,3).Value = 0
$ws.Cells.Item(55,4).Value = 16
$ws.Cells.Item(55,5).Value = 608
$ws.Cells.Item(55,6).Value = 21
$ws.Cells.Item(55,7).Value = 0
$ws.Cells.Item(55,8).Value = 11

# Row 56: Croacia
$ws.Cells.Item(56,1).Value = "Croacia"
$ws.Cells.Item(56,2).Value = 635
$ws.Cells.Item(56,3).Value = 49
$ws.Cells.Item(56,4).Value = 45
$ws.Cells.Item(56,5).Value = 586
$ws.Cells.Item(56,6).Value = 14
$ws.Cells.Item(56,7).Value = 1
$ws.Cells.Item(56,8).Value = 4

# Row 99: Sri Lanka
$ws.Cells.Item(99,1).Value = "Sri Lanka"
$ws.Cells.Item(99,2).Value = 113
$ws.Cells.Item(99,3).Value = 7
$ws.Cells.Item(99,4).Value = 9
$ws.Cells.Item(99,5).Value = 103
$ws.Cells.Item(99,6).Value = 5
$ws.Cells.Item(99,7).Value = 1
$ws.Cells.Item(99,8).Value = 1

# Row 100: Venezuela
$ws.Cells.Item(100,1).Value = "Venezuela"
$ws.Cells.Item(100,2).Value = 113
$ws.Cells.Item(100,3).Value = 0
$ws.Cells.Item(100,4).Value = 39
$ws.Cells.Item(100,5).Value = 72
$ws.Cells.Item(100,6).Value = 2
$ws.Cells.Item(100,7).Value = 0
$ws.Cells.Item(100,8).Value = 2

# Row 101: Afganistan
$ws.Cells.Item(101,1).Value = "Afganistan"
$ws.Cells.Item(101,2).Value = 110
$ws.Cells.Item(101,3).Value = 0
$ws.Cells.Item(101,4).Value = 2
$ws.Cells.Item(101,5).Value = 104
$ws.Cells.Item(101,6).Value = 0
$ws.Cells.Item(101,7).Value = 0
$ws.Cells.Item(101,8).Value = 4

# Row 159: Mozambique
$ws.Cells.Item(159,1).Value = "Mozambique"
$ws.Cells.Item(159,2).Value = 8
$ws.Cells.Item(159,3).Value = 1
$ws.Cells.Item(159,4).Value = 0
$ws.Cells.Item(159,5).Value = 8
$ws.Cells.Item(159,6).Value = 0
$ws.Cells.Item(159,7).Value = 0
$ws.Cells.Item(159,8).Value = 0

# Row 163: Guinea
$ws.Cells.Item(163,1).Value = "Guinea"
$ws.Cells.Item(163,2).Value = 8
$ws.Cells.Item(163,3).Value = 0
$ws.Cells.Item(163,4).Value = 0
$ws.Cells.Item(163,5).Value = 8
$ws.Cells.Item(163,6).Value = 0
$ws.Cells.Item(163,7).Value = 0
$ws.Cells.Item(163,8).Value = 0

# Row 164: Islas Caimanes
$ws.Cells.Item(164,1).Value = "Islas Caimanes"
$ws.Cells.Item(164,2).Value = 8
$ws.Cells.Item(164,3).Value = 0
$ws.Cells.Item(164,4).Value = 0
$ws.Cells.Item(164,5).Value = 7
$ws.Cells.Item(164,6).Value = 0
$ws.Cells.Item(164,7).Value = 0
$ws.Cells.Item(164,8).Value = 1

# Row 165: Namibia
$ws.Cells.Item(165,1).Value = "Namibia"
$ws.Cells.Item(165,2).Value = 8
$ws.Cells.Item(165,3).Value = 0
$ws.Cells.Item(165,4).Value = 2
$ws.Cells.Item(165,5).Value = 6
$ws.Cells.Item(165,6).Value = 0
$ws.Cells.Item(165,7).Value = 0
$ws.Cells.Item(165,8).Value = 0

# Row 166: Curazao
$ws.Cells.Item(166,1).Value = "Curazao"
$ws.Cells.Item(166,2).Value = 8
$ws.Cells.Item(166,3).Value = 0
$ws.Cells.Item(166,4).Value = 2
$ws.Cells.Item(166,5).Value = 5
$ws.Cells.Item(166,6).Value = 0
$ws.Cells.Item(166,7).Value = 0
$ws.Cells.Item(166,8).Value = 1

# Row 167: Seychelles
$ws.Cells.Item(167,1).Value = "Seychelles"
$ws.Cells.Item(167,2).Value = 7
$ws.Cells.Item(167,3).Value = 0
$ws.Cells.Item(167,4).Value = 0
$ws.Cells.Item(167,5).Value = 7
$ws.Cells.Item(167,6).Value = 0
$ws.Cells.Item(167,7).Value = 0
$ws.Cells.Item(167,8).Value = 0

# Row 169: Granada
$ws.Cells.Item(169,1).Value = "Granada"
$ws.Cells.Item(169,2).Value = 7
$ws.Cells.Item(169,3).Value = 0
$ws.Cells.Item(169,4).Value = 0
$ws.Cells.Item(169,5).Value = 7
$ws.Cells.Item(169,6).Value = 0
$ws.Cells.Item(169,7).Value = 0
$ws.Cells.Item(169,8).Value = 0
